# Correct typo in codebook ("1980 - 2010" -> "1981 - 2010") for the four
# Anomaly description cells, and nudge the saved selection / window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in the four "Anomaly or deviation..." description cells.
$ws.Range("F7").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 1+2  observation point (0-10°South)(90°West-80°West) "
$ws.Range("F9").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 3 observation point (5°North-5°South)(150°West-90°West)"
$ws.Range("F11").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 3.4 observation point (5°North-5°South)(170-120°West"
$ws.Range("F13").Value = "Anomaly or deviation from the 1981 - 2010 average for sea surface temperature taken at the Niño 4 observation point (5°North-5°South)"

# Update the saved selection to reflect where the cursor was left after editing.
$ws.Range("I16").Select() | Out-Null

# Nudge the saved window position (yWindow 480 -> 460) to match the author's
# view at save time.
$excel.ActiveWindow.Top = 460
